$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.382.19'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.586.49'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '2.602.28'
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.160'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.19%  '
$ws.Range("D14").Value = '3.042.44'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '59.374.51'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.64%  '
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("D18").Value = '2.593.32'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.481'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.19%  '
$ws.Range("D29").Value = '0.0₃0768'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.905'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.842'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '289.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '135.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0974'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.599'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0532'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0235'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("D50").Value = '1.967.36'
$ws.Range("E50").Value = '  +2.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.42%  '
